$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "disponible" column header (F1), bold font matching the other header cells
$ws.Range("F1").Value = "disponible"
$ws.Range("F1").Font.Bold = $true

# Fill F2:F32 with availability flag = 1
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 6).Value = 1
}

# Clear the lingering selection left in the sheet view
$null = $ws.Range("A1").Select()
